# Add a new column BB (one quarter ahead of BA) to the forecast sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date serial in row 1 (same style as the rest of the header row).
$ws.Range("BB1").Value = 45986
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null

# Copy the forecast column values from BA into the new BB column, row by row.
$bbValues = @{
    3  = -3.560752169208581
    4  = 1.224484594823672
    5  = 0.6212498672564903
    6  = 0.951852872712089
    7  = -0.3608752035976437
    8  = 0.09627146709163537
    9  = 0.1477266864992943
    10 = -0.4279125887877044
    11 = -0.002674352087272958
    12 = 0.3477863758372779
    13 = -0.8261807291073398
    14 = -1.099040380746541
    15 = 1.197694531567151
    16 = -0.7498286166554458
    17 = 0.3439499888177044
    18 = 0.2473045135454655
    19 = -2.06674933094535
    20 = -1.12081074591468
    21 = -1.194610791899986
}

foreach ($row in $bbValues.Keys) {
    $ws.Cells.Item($row, 54).Value = $bbValues[$row]
}
